# chore: update Sheets via scheduled runner
# Refresh cached market-board price/profit figures (currentAveragePrice*,
# LevePrice*/LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 303.1111
$ws.Range("I12").Value = 354.83334
$ws.Range("J12").Value = 199.66667
$ws.Range("K12").Value = 354.83334
$ws.Range("L12").Value = 199.66667
$ws.Range("M12").Value = -184.83334
$ws.Range("N12").Value = -539.6666700000001
$ws.Range("H17").Value = 1250449.1
$ws.Range("J17").Value = 1250449.1
$ws.Range("L17").Value = 3751347.3
$ws.Range("N17").Value = -3751683.3
$ws.Range("H32").Value = 5170.091
$ws.Range("I32").Value = 1833
$ws.Range("J32").Value = 6421.5
$ws.Range("K32").Value = 1833
$ws.Range("L32").Value = 6421.5
$ws.Range("M32").Value = -1507
$ws.Range("N32").Value = -7073.5
$ws.Range("H38").Value = 4280.8076
$ws.Range("I38").Value = 2770.5
$ws.Range("J38").Value = 5575.357
$ws.Range("K38").Value = 8311.5
$ws.Range("L38").Value = 16726.071
$ws.Range("M38").Value = -7939.5
$ws.Range("N38").Value = -17470.071
$ws.Range("H55").Value = 183.71428
$ws.Range("I55").Value = 180
$ws.Range("J55").Value = 186.5
$ws.Range("K55").Value = 180
$ws.Range("L55").Value = 186.5
$ws.Range("M55").Value = 34
$ws.Range("N55").Value = -614.5
$ws.Range("H80").Value = 1795.12
$ws.Range("I80").Value = 1309.625
$ws.Range("J80").Value = 2023.5883
$ws.Range("K80").Value = 3928.875
$ws.Range("L80").Value = 6070.7649
$ws.Range("M80").Value = -2930.875
$ws.Range("N80").Value = -8066.7649
$ws.Range("H83").Value = 1795.12
$ws.Range("I83").Value = 1309.625
$ws.Range("J83").Value = 2023.5883
$ws.Range("K83").Value = 11786.625
$ws.Range("L83").Value = 18212.2947
$ws.Range("M83").Value = -6794.625
$ws.Range("N83").Value = -28196.2947

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2130.75
$ws.Range("I61").Value = 2110.1428
$ws.Range("K61").Value = 2110.1428
$ws.Range("M61").Value = -1898.1428
$ws.Range("H74").Value = 35374.543
$ws.Range("I74").Value = 36296.91
$ws.Range("J74").Value = 4014
$ws.Range("K74").Value = 36296.91
$ws.Range("L74").Value = 4014
$ws.Range("M74").Value = -35422.91
$ws.Range("N74").Value = -5762
$ws.Range("H77").Value = 35374.543
$ws.Range("I77").Value = 36296.91
$ws.Range("J77").Value = 4014
$ws.Range("K77").Value = 181484.55
$ws.Range("L77").Value = 20070
$ws.Range("M77").Value = -177116.55
$ws.Range("N77").Value = -28806
$ws.Range("H88").Value = 8016.7856
$ws.Range("I88").Value = 11611.889
$ws.Range("J88").Value = 1545.6
$ws.Range("K88").Value = 11611.889
$ws.Range("L88").Value = 1545.6
$ws.Range("M88").Value = -11205.889
$ws.Range("N88").Value = -2357.6
$ws.Range("H91").Value = 8016.7856
$ws.Range("I91").Value = 11611.889
$ws.Range("J91").Value = 1545.6
$ws.Range("K91").Value = 11611.889
$ws.Range("L91").Value = 1545.6
$ws.Range("M91").Value = -10207.889
$ws.Range("N91").Value = -4353.6
$ws.Range("H136").Value = 2130.75
$ws.Range("I136").Value = 2110.1428
$ws.Range("K136").Value = 6330.428400000001
$ws.Range("M136").Value = -3780.428400000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 577
$ws.Range("J80").Value = 436.8889
$ws.Range("L80").Value = 436.8889
$ws.Range("N80").Value = -2432.8889
$ws.Range("H83").Value = 577
$ws.Range("J83").Value = 436.8889
$ws.Range("L83").Value = 2184.4445
$ws.Range("N83").Value = -12168.4445
$ws.Range("H86").Value = 26457.334
$ws.Range("I86").Value = 11548.8
$ws.Range("K86").Value = 11548.8
$ws.Range("M86").Value = -10425.8
$ws.Range("H89").Value = 26457.334
$ws.Range("I89").Value = 11548.8
$ws.Range("K89").Value = 57744
$ws.Range("M89").Value = -52128
$ws.Range("H99").Value = 2370.7334
$ws.Range("I99").Value = 2182.9285
$ws.Range("K99").Value = 2182.9285
$ws.Range("M99").Value = -684.9285
$ws.Range("H105").Value = 7755.5713
$ws.Range("I105").Value = 12981.6
$ws.Range("J105").Value = 5665.16
$ws.Range("K105").Value = 12981.6
$ws.Range("L105").Value = 5665.16
$ws.Range("M105").Value = -11234.6
$ws.Range("N105").Value = -9159.16
$ws.Range("H107").Value = 1757.5428
$ws.Range("I107").Value = 1634.9231
$ws.Range("J107").Value = 2111.7778
$ws.Range("K107").Value = 1634.9231
$ws.Range("L107").Value = 2111.7778
$ws.Range("M107").Value = 285.0769
$ws.Range("N107").Value = -5951.7778
$ws.Range("H132").Value = 118993
$ws.Range("J132").Value = 118993
$ws.Range("L132").Value = 118993
$ws.Range("N132").Value = -129113
$ws.Range("H134").Value = 3549.7856
$ws.Range("I134").Value = 2425.818
$ws.Range("K134").Value = 7277.454000000001
$ws.Range("M134").Value = -4742.454000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2419.5557
$ws.Range("I16").Value = 2353.7144
$ws.Range("J16").Value = 2650
$ws.Range("K16").Value = 2353.7144
$ws.Range("L16").Value = 2650
$ws.Range("M16").Value = -2066.7144
$ws.Range("N16").Value = -3224
$ws.Range("H31").Value = 3943.1765
$ws.Range("I31").Value = 2186.2222
$ws.Range("K31").Value = 2186.2222
$ws.Range("M31").Value = -1891.2222
$ws.Range("H34").Value = 3943.1765
$ws.Range("I34").Value = 2186.2222
$ws.Range("K34").Value = 2186.2222
$ws.Range("M34").Value = -1984.2222
$ws.Range("H107").Value = 355
$ws.Range("I107").Value = 11
$ws.Range("J107").Value = 699
$ws.Range("K107").Value = 11
$ws.Range("L107").Value = 699
$ws.Range("M107").Value = 1909
$ws.Range("N107").Value = -4539
$ws.Range("H113").Value = 2419.5557
$ws.Range("I113").Value = 2353.7144
$ws.Range("J113").Value = 2650
$ws.Range("K113").Value = 2353.7144
$ws.Range("L113").Value = 2650
$ws.Range("M113").Value = -183.7143999999998
$ws.Range("N113").Value = -6990

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5666.6665
$ws.Range("I68").Value = 15000
$ws.Range("K68").Value = 45000
$ws.Range("M68").Value = -44189
$ws.Range("H71").Value = 5666.6665
$ws.Range("I71").Value = 15000
$ws.Range("K71").Value = 135000
$ws.Range("M71").Value = -130944
$ws.Range("H132").Value = 4367.7144
$ws.Range("I132").Value = 1314.8
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 11833.2
$ws.Range("L132").Value = 108000
$ws.Range("M132").Value = -9303.199999999999
$ws.Range("N132").Value = -113060

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 29685832
$ws.Range("I11").Value = 44513748
$ws.Range("K11").Value = 44513748
$ws.Range("M11").Value = -44513609
$ws.Range("H12").Value = 5000
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5280
$ws.Range("H97").Value = 1904.091
$ws.Range("I97").Value = 906.8
$ws.Range("J97").Value = 2735.1667
$ws.Range("K97").Value = 906.8
$ws.Range("L97").Value = 2735.1667
$ws.Range("M97").Value = -410.8
$ws.Range("N97").Value = -3727.1667

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7581.5
$ws.Range("I7").Value = 6999.6665
$ws.Range("J7").Value = 8163.3335
$ws.Range("K7").Value = 6999.6665
$ws.Range("L7").Value = 8163.3335
$ws.Range("M7").Value = -6887.6665
$ws.Range("N7").Value = -8387.333500000001
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("H43").Value = 20836.666
$ws.Range("I43").Value = 21500
$ws.Range("K43").Value = 21500
$ws.Range("M43").Value = -21307
$ws.Range("H46").Value = 2999.9048
$ws.Range("J46").Value = 4499.8335
$ws.Range("L46").Value = 4499.8335
$ws.Range("N46").Value = -4875.8335
$ws.Range("H122").Value = 4209.5
$ws.Range("I122").Value = 3974.348
$ws.Range("J122").Value = 4810.4443
$ws.Range("K122").Value = 11923.044
$ws.Range("L122").Value = 14431.3329
$ws.Range("M122").Value = -9473.044
$ws.Range("N122").Value = -19331.3329
$ws.Range("H126").Value = 7581.5
$ws.Range("I126").Value = 6999.6665
$ws.Range("J126").Value = 8163.3335
$ws.Range("K126").Value = 20998.9995
$ws.Range("L126").Value = 24490.0005
$ws.Range("M126").Value = -18528.9995
$ws.Range("N126").Value = -29430.0005

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2003302.2
$ws.Range("I23").Value = 2504075
$ws.Range("J23").Value = 211
$ws.Range("K23").Value = 2504075
$ws.Range("L23").Value = 211
$ws.Range("M23").Value = -2503846
$ws.Range("N23").Value = -669
$ws.Range("H132").Value = 5124.0625
$ws.Range("I132").Value = 5065.6665
$ws.Range("K132").Value = 15196.9995
$ws.Range("M132").Value = -12666.9995
$ws.Range("H133").Value = 79578.25
$ws.Range("J133").Value = 79578.25
$ws.Range("L133").Value = 79578.25
$ws.Range("N133").Value = -89698.25
